$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

$xlPasteFormats = -4122

# --- Row 4 (new): reading collections/linq note (entered before the row-3
#     wording tweak so the shared-string table picks up the same ordering
#     the authored workbook has) ---
$ws.Range("A2").Copy()
$ws.Range("A4").PasteSpecial($xlPasteFormats)
$ws.Range("A4").Value = 43809
$ws.Range("B4").Value = "Read into Collections, Linq, Lambda expressions"
$ws.Range("C4").Value = 3

# --- Row 3: fix the wording and shrink the row height ---
$ws.Range("B3").Value = 'Created "Consignment Shop" App. Research on Properties, Interfaces and Patterns.'
$ws.Rows.Item(3).RowHeight = 30

# --- Row 5 (new): code review note ---
$ws.Range("A2").Copy()
$ws.Range("A5").PasteSpecial($xlPasteFormats)
$ws.Range("A5").Value = 43810
$ws.Range("B2").Copy()
$ws.Range("B5").PasteSpecial($xlPasteFormats)
$ws.Range("B5").Value = 'Review of "Consignment Shop" code, Intial code refactoring, Introduction to basic principles.'
$ws.Range("C5").Value = 2
$ws.Rows.Item(5).RowHeight = 45

$excel.CutCopyMode = $false

# --- column widths (target OOXML widths are 13 and 44.140625 chars) ---
$ws.Columns.Item(1).ColumnWidth = 12.166666666666668
$ws.Columns.Item(2).ColumnWidth = 43.333333333333336

# --- selection moves to B6 after the new rows ---
$ws.Range("B6").Select()
